$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update rows 2-10 with the latest Banshee relay settings.
for ($r = 2; $r -le 10; $r++) {
    $ws.Range("F$r").Value = 13
    $ws.Range("N$r").Value = 0.7
    $ws.Range("Q$r").Value = 0

    # R column takes on the 51P TOC Time Dial value already computed in column K,
    # and picks up that column's number formatting (2 decimal places).
    $kVal = $ws.Range("K$r").Value2
    $ws.Range("R$r").NumberFormat = $ws.Range("K$r").NumberFormat
    $ws.Range("R$r").Value = $kVal

    $ws.Range("S$r").Value = 1
    $ws.Range("T$r").Value = 1
    $ws.Range("U$r").Value = 0.5
    $ws.Range("V$r").Value = 1.2
}

# Update the active selection left over from editing.
$ws.Range("H12").Select()
